$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. optimization_parameters sheet content changes
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the extra repeated "value" header cells (C1:F1)
$ws.Range("C1:F1").ClearContents()

# Row 8 label: "Model" -> "production_function"
$ws.Range("A8").Value = "production_function"

# Insert a new row 9 ("L_curve") - shifts old rows 9-17 down to 10-18
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = "0.00E+00"

# Remove the old "Deletion" row (now at row 17, after the insert above)
$ws.Rows.Item(17).Delete()

# Tab/selection changes on this sheet: it becomes the active/selected tab
$ws.Activate()
$ws.Range("C1:I2").Select()

# ------------------------------------------------------------------
# 2. network_weights sheet loses the "selected tab" flag
#    (handled implicitly by activating optimization_parameters above,
#    but make sure its own selection stays where it was)
# ------------------------------------------------------------------
$wsnw = $wb.Worksheets.Item("network_weights")
$wsnw.Range("E12").Select()

# ------------------------------------------------------------------
# 3. Workbook view: shift the first visible tab / active tab by one
# ------------------------------------------------------------------
$wb.Windows.Item(1).ScrollWorkbookTabs(1)

# Re-activate optimization_parameters so it is the active sheet on open
$ws.Activate()

Write-Output "edit complete"
